$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) column names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Title-case the Spanish connector words (de/del/la/las/el/los/y) in state/municipality names
$ws.Range('B6').Value = 'Pabellón De Arteaga'
$ws.Range('B25').Value = 'Amatenango De La Frontera'
$ws.Range('B26').Value = 'Amatenango Del Valle'
$ws.Range('B29').Value = 'Bejucal De Ocampo'
$ws.Range('B31').Value = 'Benemérito De Las Américas'
$ws.Range('B38').Value = 'Comitán De Domínguez'
$ws.Range('B54').Value = 'Marqués De Comillas'
$ws.Range('B55').Value = 'Mazapa De Madero'
$ws.Range('B57').Value = 'Montecristo De Guerrero'
$ws.Range('B63').Value = 'Salto De Agua'
$ws.Range('B64').Value = 'San Cristóbal De Las Casas'
$ws.Range('B84').Value = 'Valle De Zaragoza'
$ws.Range('B96').Value = 'San Juan De Sabinas'
$ws.Range('A105').Value = 'Ciudad De México'
$ws.Range('B109').Value = 'Cuajimalpa De Morelos'
$ws.Range('B128').Value = 'Nombre De Dios'
$ws.Range('B130').Value = 'Pánuco De Coronado'
$ws.Range('A143').Value = 'Estado De México'
$ws.Range('B143').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B145').Value = 'Almoloya De Juárez'
$ws.Range('B150').Value = 'Atizapán De Zaragoza'
$ws.Range('B155').Value = 'Coacalco De Berriozábal'
$ws.Range('B160').Value = 'Ecatepec De Morelos'
$ws.Range('B165').Value = 'Ixtapan Del Oro'
$ws.Range('B172').Value = 'Naucalpan De Juárez'
$ws.Range('B178').Value = 'San Felipe Del Progreso'
$ws.Range('B180').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B188').Value = 'Tenango Del Valle'
$ws.Range('B194').Value = 'Tlalnepantla De Baz'
$ws.Range('B198').Value = 'Villa De Allende'
$ws.Range('B199').Value = 'Villa Del Carbón'
$ws.Range('B208').Value = 'Apaseo El Alto'
$ws.Range('B209').Value = 'Apaseo El Grande'
$ws.Range('B213').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B217').Value = 'Jaral Del Progreso'
$ws.Range('B227').Value = 'San Diego De La Unión'
$ws.Range('B229').Value = 'San Francisco Del Rincón'
$ws.Range('B231').Value = 'San Luis De La Paz'
$ws.Range('B232').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B233').Value = 'Silao De La Victoria'
$ws.Range('B235').Value = 'Valle De Santiago'
$ws.Range('B239').Value = 'Acapulco De Juárez'
$ws.Range('B241').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B242').Value = 'Alcozauca De Guerrero'
$ws.Range('B245').Value = 'Atenango Del Río'
$ws.Range('B246').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B248').Value = 'Atoyac De Álvarez'
$ws.Range('B249').Value = 'Ayutla De Los Libres'
$ws.Range('B251').Value = 'Buenavista De Cuéllar'
$ws.Range('B252').Value = 'Chilapa De Álvarez'
$ws.Range('B253').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B256').Value = 'Coyuca De Benítez'
$ws.Range('B257').Value = 'Coyuca De Catalán'
$ws.Range('B260').Value = 'Cuetzala Del Progreso'
$ws.Range('B261').Value = 'Cutzamala De Pinzón'
$ws.Range('B265').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B266').Value = 'Iguala De La Independencia'
$ws.Range('B268').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B269').Value = 'Zihuatanejo De Azueta'
$ws.Range('B271').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B285').Value = 'Taxco De Alarcón'
$ws.Range('B287').Value = 'Técpan De Galeana'
$ws.Range('B289').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B290').Value = 'Tixtla De Guerrero'
$ws.Range('B293').Value = 'Tlapa De Comonfort'
$ws.Range('B302').Value = 'Atotonilco De Tula'
$ws.Range('B303').Value = 'Atotonilco El Grande'
$ws.Range('B306').Value = 'Cuautepec De Hinojosa'
$ws.Range('B309').Value = 'Huejutla De Reyes'
$ws.Range('B312').Value = 'Jacala De Ledezma'
$ws.Range('B317').Value = 'Mixquiahuala De Juárez'
$ws.Range('B318').Value = 'Molango De Escamilla'
$ws.Range('B320').Value = 'Pachuca De Soto'
$ws.Range('B323').Value = 'Progreso De Obregón'
$ws.Range('B326').Value = 'Santiago De Anaya'
$ws.Range('B329').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B330').Value = 'Tezontepec De Aldama'
$ws.Range('B335').Value = 'Tula De Allende'
$ws.Range('B336').Value = 'Tulancingo De Bravo'
$ws.Range('B337').Value = 'Zacualtipán De Ángeles'
$ws.Range('B347').Value = 'Encarnación De Díaz'
$ws.Range('B349').Value = 'Huejuquilla El Alto'
$ws.Range('B350').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B352').Value = 'Jilotlán De Los Dolores'
$ws.Range('B355').Value = 'Lagos De Moreno'
$ws.Range('B359').Value = 'Ojuelos De Jalisco'
$ws.Range('B362').Value = 'San Juan De Los Lagos'
$ws.Range('B363').Value = 'San Juanito De Escobedo'
$ws.Range('B366').Value = 'San Miguel El Alto'
$ws.Range('B367').Value = 'Santa María Del Oro'
$ws.Range('B369').Value = 'Tamazula De Gordiano'
$ws.Range('B371').Value = 'Tepatitlán De Morelos'
$ws.Range('B377').Value = 'Unión De Tula'
$ws.Range('B380').Value = 'Yahualica De González Gallo'
$ws.Range('B399').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B439').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B466').Value = 'Tlaltizapán De Zapata'
$ws.Range('B484').Value = 'Mier Y Noriega'
$ws.Range('B487').Value = 'San Nicolás De Los Garza'
$ws.Range('B489').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B495').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B496').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B497').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B498').Value = 'Ixtlán De Juárez'
$ws.Range('B499').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B503').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B505').Value = 'Oaxaca De Juárez'
$ws.Range('B506').Value = 'Ocotlán De Morelos'
$ws.Range('B508').Value = 'Putla Villa De Guerrero'
$ws.Range('B516').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B538').Value = 'San Miguel Del Puerto'
$ws.Range('B546').Value = 'San Pablo Villa De Mitla'
$ws.Range('B547').Value = 'San Pedro El Alto'
$ws.Range('B589').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B590').Value = 'Tataltepec De Valdés'
$ws.Range('B591').Value = 'Teococuilco De Marcos Pérez'
$ws.Range('B592').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B593').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B594').Value = 'Tlacolula De Matamoros'
$ws.Range('B595').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B602').Value = 'Ayotoxco De Guerrero'
$ws.Range('B604').Value = 'Chalchicomula De Sesma'
$ws.Range('B623').Value = 'Izúcar De Matamoros'
$ws.Range('B638').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B640').Value = 'San Salvador El Seco'
$ws.Range('B643').Value = 'Tepango De Rodríguez'
$ws.Range('B646').Value = 'Tepexi De Rodríguez'
$ws.Range('B647').Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range('B648').Value = 'Tetela De Ocampo'
$ws.Range('B651').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B664').Value = 'Amealco De Bonfil'
$ws.Range('B665').Value = 'Cadereyta De Montes'
$ws.Range('B669').Value = 'Jalpan De Serra'
$ws.Range('B670').Value = 'Landa De Matamoros'
$ws.Range('B671').Value = 'Pinal De Amoles'
$ws.Range('B674').Value = 'San Juan Del Río'
$ws.Range('B684').Value = 'Ciudad Del Maíz'
$ws.Range('B695').Value = 'San Ciro De Acosta'
$ws.Range('B698').Value = 'Santa María Del Río'
$ws.Range('B704').Value = 'Tanquián De Escobedo'
$ws.Range('B708').Value = 'Villa De Arista'
$ws.Range('B709').Value = 'Villa De Arriaga'
$ws.Range('B710').Value = 'Villa De La Paz'
$ws.Range('B711').Value = 'Villa De Ramos'
$ws.Range('B712').Value = 'Villa De Reyes'
$ws.Range('B738').Value = 'Jalpa De Méndez'
$ws.Range('B759').Value = 'Soto La Marina'
$ws.Range('B766').Value = 'Apetatitlán De Antonio Carvajal'
$ws.Range('B771').Value = 'San Pablo Del Monte'
$ws.Range('B774').Value = 'Tetla De La Solidaridad'
$ws.Range('B783').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B786').Value = 'Amatlán De Los Reyes'
$ws.Range('B796').Value = 'Castillo De Teayo'
$ws.Range('B807').Value = 'Cosamaloapan De Carpio'
$ws.Range('B808').Value = 'Cosautlán De Carvajal'
$ws.Range('B824').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B825').Value = 'Ixhuatlán De Madero'
$ws.Range('B826').Value = 'Ixhuatlán Del Café'
$ws.Range('B833').Value = 'Juchique De Ferrer'
$ws.Range('B837').Value = 'Lerdo De Tejada'
$ws.Range('B840').Value = 'Martínez De La Torre'
$ws.Range('B844').Value = 'Mixtla De Altamirano'
$ws.Range('B852').Value = 'Paso Del Macho'
$ws.Range('B855').Value = 'Poza Rica De Hidalgo'
$ws.Range('B860').Value = 'Sayula De Alemán'
$ws.Range('B862').Value = 'Soledad De Doblado'
$ws.Range('B881').Value = 'Vega De Alatorre'
$ws.Range('B889').Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range('B891').Value = 'Dzilam De Bravo'
$ws.Range('B898').Value = 'El Plateado De Joaquín Amaro'

# Remove trailing metadata/footer rows (916-920)
$ws.Range('A916:A920').EntireRow.Delete()
